$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly scoreboard rows (week 7) being appended below the existing data.
$newRows = @(
    @{ A = "Matt";     B = 45496; C = "Run";  D = 41; E = 4.77; F = 276; G = 0;  H = 15; I = 11; J = 12; K = 0; L = "Agile Antelope";    M = 7 },
    @{ A = "Matt";     B = 45496; C = "Walk"; D = 3;  E = 0.13; F = 0;   G = 3;  H = 0;  I = 0;  J = 0;  K = 0; L = "Agile Antelope";    M = 7 },
    @{ A = "Jeremiah"; B = 45496; C = "Run";  D = 20; E = 2.24; F = 104; G = 0;  H = 11; I = 7;  J = 0;  K = 0; L = "Sauntering Hippo";  M = 7 },
    @{ A = "Steven";   B = 45496; C = "Walk"; D = 41; E = 2.16; F = 89;  G = 41; H = 0;  I = 0;  J = 0;  K = 0; L = "Brave Leopard";     M = 7 },
    @{ A = "Steven";   B = 45496; C = "Walk"; D = 31; E = 1.44; F = 49;  G = 31; H = 0;  I = 0;  J = 0;  K = 0; L = "Brave Leopard";     M = 7 },
    @{ A = "Steven";   B = 45497; C = "Walk"; D = 27; E = 1.33; F = 66;  G = 27; H = 0;  I = 0;  J = 0;  K = 0; L = "Brave Leopard";     M = 7 }
)

$startRow = 264
$endRow = $startRow + $newRows.Count - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $row.B
    $ws.Cells.Item($r, 3).Value2 = $row.C
    $ws.Cells.Item($r, 4).Value2 = $row.D
    $ws.Cells.Item($r, 5).Value2 = $row.E
    $ws.Cells.Item($r, 6).Value2 = $row.F
    $ws.Cells.Item($r, 7).Value2 = $row.G
    $ws.Cells.Item($r, 8).Value2 = $row.H
    $ws.Cells.Item($r, 9).Value2 = $row.I
    $ws.Cells.Item($r, 10).Value2 = $row.J
    $ws.Cells.Item($r, 11).Value2 = $row.K
    $ws.Cells.Item($r, 12).Value2 = $row.L
    $ws.Cells.Item($r, 13).Value2 = $row.M
}

# The date column (B) uses the same short-date display format as the rest of
# the table; copy that formatting down onto the newly added rows so the new
# cells share the existing style record instead of creating a new one.
$ws.Cells.Item($startRow - 1, 2).Copy()
$ws.Range("B$startRow`:B$endRow").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Application.CutCopyMode = $false

# Move the active selection to reflect where the user would type the next entry.
$ws.Range("D270").Select()
